# "First Commit through VS Code - startwork branch"
#
# Adds a new note after the "Committed without the "-m" message..."
# paragraph:
#   (blank paragraph)
#   Git checkout [Wingdings right-arrow symbol] Switches branches
#   (blank paragraph) x3
#
# ...directly ahead of the two pre-existing trailing blank paragraphs.

$d = $word.ActiveDocument

# Locate the paragraph that ends the "Committed without ..." note - the
# anchor point after which the new content is inserted.
$anchor = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Committed without*") {
        $anchor = $para
    }
}

if ($anchor -eq $null) {
    throw "Could not find the 'Committed without ...' paragraph to anchor the insertion."
}

# Create a fresh paragraph right after the anchor; its Range is then
# replaced (via WordprocessingML) with the five paragraphs we need. Doing
# the split first (InsertParagraphAfter) keeps the anchor paragraph intact
# instead of letting the OOXML insertion merge into it.
$null = $anchor.Range.InsertParagraphAfter()
$newPara = $anchor.Next()

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$blankP = "<w:p $wNs/>"
$gitCheckoutP = "<w:p $wNs>" + `
    '<w:r><w:t xml:space="preserve">Git checkout </w:t></w:r>' + `
    '<w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t>Switches branches</w:t></w:r>' + `
    "</w:p>"

$fragment = $blankP + $gitCheckoutP + $blankP + $blankP + $blankP

$null = $newPara.Range.InsertXML($fragment)

Write-Output ("Paragraphs after edit: " + $d.Paragraphs.Count)
